$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 243, shifting existing rows 243:283 down to 245:285
$ws.Range("A243:A244").EntireRow.Insert()

# Row 243 - new record
$ws.Range("A243").Value = 1
$ws.Range("B243").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C243").Value = "Arica y Parinacota"
$ws.Range("D243").Value = 45204
$ws.Range("E243").Value = 15
$ws.Range("F243").Value = 100114001
$ws.Range("G243").Value = "Papa"
$ws.Range("H243").Value = "Asterix"
$ws.Range("I243").Value = "1a (cosecha)"
$ws.Range("J243").Value = 550
$ws.Range("K243").Value = 32000
$ws.Range("L243").Value = 33000
$ws.Range("M243").Value = 32455
$ws.Range("N243").Value = "`$/saco 25 kilos"
$ws.Range("O243").Value = "Región Metropolitana"
$ws.Range("P243").Value = 1298
$ws.Range("Q243").Value = 25
$ws.Range("R243").Value = "Hortaliza"

# Row 244 - new record
$ws.Range("A244").Value = 1
$ws.Range("B244").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C244").Value = "Arica y Parinacota"
$ws.Range("D244").Value = 45204
$ws.Range("E244").Value = 15
$ws.Range("F244").Value = 100114001
$ws.Range("G244").Value = "Papa"
$ws.Range("H244").Value = "Cardinal"
$ws.Range("I244").Value = "1a (cosecha)"
$ws.Range("J244").Value = 750
$ws.Range("K244").Value = 33000
$ws.Range("L244").Value = 34000
$ws.Range("M244").Value = 33533
$ws.Range("N244").Value = "`$/saco 25 kilos"
$ws.Range("O244").Value = "Región de Coquimbo"
$ws.Range("P244").Value = 1341
$ws.Range("Q244").Value = 25
$ws.Range("R244").Value = "Hortaliza"
